# Update the Wnt5a-Fzd6 LR-pairs sheet with newly-computed TPM-based values.
# The sending/target cluster set changed from {FAPs -> ECs,FAPs,MuSCs,Neutrophils}
# to {ECs,FAPs} x {ECs,FAPs,MuSCs} (Neutrophils target dropped, ECs sender added),
# and all of the expression/specificity metrics were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $vals)
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Sending cluster | Ligand symbol | Receptor symbol | Target cluster | ... metrics ...
Set-Row 2 @("ECs", "Wnt5a", "Fzd6", "ECs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 2, 1, 11.4621665, 22.924333, 0.9441106217120846, 0.9320562266826696, 1.4398200448975, 5.75928017959, 0.0249951273549395, 0.01659712835611562)
Set-Row 3 @("ECs", "Wnt5a", "Fzd6", "FAPs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 2, 0.6666666666666666, 0.3140343333333334, 0.942103, 0.02586624000639611, 0.03830397016682768, 0.03944742278166667, 0.23668453669, 0.0006848031875553534, 0.0006820789252922471)
Set-Row 4 @("ECs", "Wnt5a", "Fzd6", "MuSCs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 2, 1, 0.364502, 0.729004, 0.03002313828151931, 0.02963980315050269, 0.04578691873, 0.18314767492, 0.000794856182828103, 0.0005277960741593533)
Set-Row 5 @("FAPs", "Wnt5a", "Fzd6", "ECs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 2, 1, 11.4621665, 22.924333, 0.9441106217120846, 0.9320562266826696, 52.94475573415201, 317.668534404912, 0.9191154943571451, 0.915459098326554)
Set-Row 6 @("FAPs", "Wnt5a", "Fzd6", "FAPs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 2, 0.6666666666666666, 0.3140343333333334, 0.942103, 0.02586624000639611, 0.03830397016682768, 1.450552220688, 13.054969986192, 0.02518143681884076, 0.03762189124153543)
Set-Row 7 @("FAPs", "Wnt5a", "Fzd6", "MuSCs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 2, 1, 0.364502, 0.729004, 0.03002313828151931, 0.02963980315050269, 1.683666814176, 10.102000885056, 0.02922828209869121, 0.02911200707634334)
